# Update the AccountKeyword value in B2 from "Tea" to "Snacks"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountLevelKeywords")

$ws.Range("B2").Value = "Snacks"

# Move/leave the active selection on the edited cell (B2)
$ws.Range("B2").Select()
